# Add new daily-track rows (13th entry, dated 17-Mar-2022) to the MAR-22
# sheet, mirroring the style/layout already used by the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAR-22")   # the "MAR-22" tab is the one being edited
$ws.Activate()

# --- Row 28 -----------------------------------------------------------
# No / Date / Application columns for the new day entry.
$ws.Cells.Item(28, 1).Value = 13

$ws.Cells.Item(27, 2).Copy() | Out-Null
$ws.Cells.Item(28, 2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats (date fmt)
$ws.Cells.Item(28, 2).Value = 44637

$ws.Cells.Item(28, 3).Value = "RPA GSS"

$ws.Cells.Item(26, 4).Copy() | Out-Null
$ws.Cells.Item(28, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats (wrap text)
$ws.Cells.Item(28, 4).Value = "1. Correction Received for the Service Order Pending Task to implement the reason with aging status in the filename during file moving to backup folder and also to write all uploaded filenames at log file and it has been completed, tested and it is`r`nrunning smoothly"

$ws.Cells.Item(27, 5).Copy() | Out-Null
$ws.Cells.Item(28, 5).PasteSpecial(-4122) | Out-Null   # xlPasteFormats (percent fmt)
$ws.Cells.Item(28, 5).Value = 1

$ws.Cells.Item(28, 6).Value = "Completed"

# Row 28 wraps across four lines once the comment is filled in.
$ws.Rows.Item(28).RowHeight = 57.6

# --- Row 29 -------------------------------------------------------------
$ws.Cells.Item(29, 4).Value = "2. Implementation of the public holidays is work in progress"

$ws.Cells.Item(27, 5).Copy() | Out-Null
$ws.Cells.Item(29, 5).PasteSpecial(-4122) | Out-Null   # xlPasteFormats (percent fmt)
$ws.Cells.Item(29, 5).Value = 0.2

$ws.Cells.Item(29, 6).Value = "WIP"

# --- Row 30 -------------------------------------------------------------
$ws.Cells.Item(30, 3).Value = "RPA RLOGIC"
$ws.Cells.Item(30, 4).Value = "1. Correction updated at P&L reports at MLR for the Jan22"

# --- Selection / scroll position matching the saved view ---------------
$excel.Goto($ws.Range("A19"), $true) | Out-Null
$ws.Range("D34").Select() | Out-Null
